# Update column G ("K") values on the active sheet to reflect the
# regenerated save_data (K instead of Strike#, std/mean recalculated,
# s_vals recalculated). Only the literal values in G2:G30 change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 7
    3  = 1
    4  = 3
    5  = 4
    6  = 3
    7  = 3
    8  = 1
    9  = 3
    10 = 2
    11 = 1
    12 = 6
    13 = 3
    14 = 4
    15 = 2
    16 = 3
    17 = 6
    18 = 3
    19 = 4
    20 = 2
    21 = 5
    22 = 1
    23 = 4
    24 = 0
    25 = 2
    26 = 6
    27 = 3
    28 = 3
    29 = 2
    30 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
